$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking price strings
# (e.g. '1.00', '584.16') are preserved exactly as text, matching the source data.
$ws.Range('D2:D51').NumberFormat = "@"

$ws.Range('D2').Value = '70.130.56'
$ws.Range('E2').Value = '  +3.23%  '
$ws.Range('D3').Value = '3.403.71'
$ws.Range('E3').Value = '  +2.08%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '584.16'
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('D6').Value = '181.99'
$ws.Range('E6').Value = '  +2.50%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').Value = '0.597'
$ws.Range('E8').Value = '  +0.87%  '
$ws.Range('E9').Value = '  +10.22%  '
$ws.Range('D10').Value = '0.595'
$ws.Range('E10').Value = '  +1.90%  '
$ws.Range('D11').Value = '48.85'
$ws.Range('E11').Value = '  +1.31%  '
$ws.Range('D12').Value = '0.0000288'
$ws.Range('E12').Value = '  +5.26%  '
$ws.Range('D13').Value = '688.89'
$ws.Range('E13').Value = '  -1.38%  '
$ws.Range('D14').Value = '8.70'
$ws.Range('E14').Value = '  +3.07%  '
$ws.Range('D15').Value = '3.956.75'
$ws.Range('E15').Value = '  +2.04%  '
$ws.Range('D16').Value = '70.090.39'
$ws.Range('E16').Value = '  +3.12%  '
$ws.Range('D17').Value = '3.407.51'
$ws.Range('E17').Value = '  +1.97%  '
$ws.Range('E18').Value = '  +1.23%  '
$ws.Range('D19').Value = '17.73'
$ws.Range('E19').Value = '  +1.18%  '
$ws.Range('D20').Value = '11.39'
$ws.Range('E20').Value = '  +2.05%  '
$ws.Range('D21').Value = '0.918'
$ws.Range('E21').Value = '  +2.52%  '
$ws.Range('D22').Value = '17.41'
$ws.Range('E22').Value = '  +2.88%  '
$ws.Range('D23').Value = '5.37'
$ws.Range('E23').Value = '  -0.54%  '
$ws.Range('D24').Value = '102.59'
$ws.Range('E24').Value = '  +2.37%  '
$ws.Range('D25').Value = '3.93'
$ws.Range('E25').Value = '  +0.52%  '
$ws.Range('D26').Value = '2.72'
$ws.Range('E26').Value = '  +1.12%  '
$ws.Range('D27').Value = '9.73'
$ws.Range('E27').Value = '  +2.59%  '
$ws.Range('D28').Value = '33.90'
$ws.Range('E28').Value = '  +2.68%  '
$ws.Range('D29').Value = '8.83'
$ws.Range('E29').Value = '  +3.33%  '
$ws.Range('D30').Value = '6.99'
$ws.Range('E30').Value = '  +0.29%  '
$ws.Range('D31').Value = '3.74'
$ws.Range('E31').Value = '  +10.94%  '
$ws.Range('D32').Value = '11.16'
$ws.Range('E32').Value = '  +0.98%  '
$ws.Range('D33').Value = '557.50'
$ws.Range('E33').Value = '  -3.40%  '
$ws.Range('E34').Value = '  +1.37%  '
$ws.Range('D35').Value = '58.72'
$ws.Range('E35').Value = '  +2.29%  '
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('D37').Value = '3.661.12'
$ws.Range('E37').Value = '  -2.28%  '
$ws.Range('E38').Value = '  +3.54%  '
$ws.Range('D39').Value = '35.68'
$ws.Range('E39').Value = '  +0.91%  '
$ws.Range('D40').Value = '0.0₃0736'
$ws.Range('E40').Value = '  +8.96%  '
$ws.Range('D41').Value = '3.34'
$ws.Range('E41').Value = '  +5.57%  '
$ws.Range('D42').Value = '2.75'
$ws.Range('E42').Value = '  +4.57%  '
$ws.Range('E43').Value = '  +4.66%  '
$ws.Range('E44').Value = '  +1.66%  '
$ws.Range('E45').Value = '  +2.65%  '
$ws.Range('E46').Value = '  +0.83%  '
$ws.Range('E47').Value = '  +4.47%  '
$ws.Range('D48').Value = '1.00'
$ws.Range('E48').Value = '  -0.11%  '
$ws.Range('D49').Value = '130.80'
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('D50').Value = '2.61'
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('D51').Value = '7.56'
$ws.Range('E51').Value = '  +2.54%  '
